$d = $word.ActiveDocument

# Locate the "Ver no Jupiter..." paragraph, which marks the beginning of the
# boilerplate footer block that needs to be removed (that paragraph, the
# copyright paragraph right after it, and the blank paragraph that trails
# them, right up to the page-break paragraph that closes the document).
$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Resolve the paragraph index containing the found text so we can
    # extend the deletion to cover the following two paragraphs as well.
    $startPara = $findRange.Paragraphs.Item(1)
    $startParaIndex = 0
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Start -eq $startPara.Range.Start) {
            $startParaIndex = $i
            break
        }
    }

    if ($startParaIndex -gt 0) {
        $deleteStart = $d.Paragraphs.Item($startParaIndex).Range.Start
        $deleteEnd = $d.Paragraphs.Item($startParaIndex + 2).Range.End
        $r = $d.Range($deleteStart, $deleteEnd)
        $r.Delete()
    }
}
